# New crime data collected - update the 6th Precinct weekly CompStat report
# with the values for the week covering 5/5/2025 through 5/11/2025
# (Volume 32, Number 19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper functions used to turn a numeric cell into a "text" cell that
# stores one of the two special shared placeholder strings used
# throughout this report ("0" or "***.*") while keeping the same
# "label" cell style (s=13) that is used elsewhere in the sheet for
# this kind of placeholder.
#
# A plain `.Value = "0"` assignment gets auto-detected by Excel as a
# number, and `.Value = "***.*"` (while text) would still leave the
# cell on its old (numeric) style. Copying the format from a cell that
# already uses the desired "label" style, then copying the value from
# a cell that already contains the desired placeholder text,
# reproduces the exact target state.
# ---------------------------------------------------------------------
function Set-TextZero([string]$addr) {
    $ws.Range("C23").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range("C23").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4163) | Out-Null   # xlPasteValues
}
function Set-TextStar([string]$addr) {
    $ws.Range("E23").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range("E23").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# ---------------------------------------------------------------------
# Header: volume/number and reporting week dates
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  19"
$ws.Range("C9").Value = "Report Covering the Week  5/5/2025  Through  5/11/2025"

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 32
$ws.Range("J16").Value = 56
$ws.Range("K16").Value = -42.857142857142
$ws.Range("L16").Value = -47.540983606557
$ws.Range("M16").Value = -30.434782608695
$ws.Range("N16").Value = -89.333333333333

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
Set-TextZero "C17"
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 11
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = -31.25
$ws.Range("I17").Value = 38
$ws.Range("J17").Value = 40
$ws.Range("K17").Value = -5
$ws.Range("L17").Value = -36.666666666666
$ws.Range("M17").Value = 18.75
$ws.Range("N17").Value = -55.294117647058

# ---------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 11
$ws.Range("E18").Value = -81.818181818181
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 29
$ws.Range("H18").Value = -55.172413793103
$ws.Range("I18").Value = 54
$ws.Range("J18").Value = 93
$ws.Range("K18").Value = -41.935483870967
$ws.Range("L18").Value = -52.212389380531
$ws.Range("M18").Value = -26.027397260274
$ws.Range("N18").Value = -80

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 21
$ws.Range("D19").Value = 22
$ws.Range("E19").Value = -4.545454545454
$ws.Range("F19").Value = 68
$ws.Range("G19").Value = 81
$ws.Range("H19").Value = -16.049382716049
$ws.Range("I19").Value = 319
$ws.Range("J19").Value = 355
$ws.Range("K19").Value = -10.140845070422
$ws.Range("L19").Value = -26.327944572748
$ws.Range("M19").Value = -8.333333333333
$ws.Range("N19").Value = -60.273972602739

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
Set-TextZero "D20"
Set-TextStar "E20"
$ws.Range("N20").Value = -98.168498168498

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 41
$ws.Range("E21").Value = -36.585365853658
$ws.Range("F21").Value = 106
$ws.Range("G21").Value = 140
$ws.Range("H21").Value = -24.285714285714
$ws.Range("I21").Value = 453
$ws.Range("J21").Value = 561
$ws.Range("K21").Value = -19.251336898395
$ws.Range("L21").Value = -33.577712609970
$ws.Range("M21").Value = -12.038834951456
$ws.Range("N21").Value = -73.905529953917

# ---------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------
Set-TextZero "C22"
$ws.Range("F22").Value = 3
Set-TextZero "G22"
Set-TextStar "H22"
$ws.Range("I22").Value = 20
$ws.Range("K22").Value = 25
$ws.Range("L22").Value = 11.111111111111
$ws.Range("M22").Value = -20

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = -7.692307692307
$ws.Range("F24").Value = 99
$ws.Range("G24").Value = 140
$ws.Range("H24").Value = -29.285714285714
$ws.Range("I24").Value = 513
$ws.Range("J24").Value = 612
$ws.Range("K24").Value = -16.176470588235
$ws.Range("L24").Value = -18.957345971564
$ws.Range("M24").Value = -2.099236641221

# ---------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 20
$ws.Range("E25").Value = 10
$ws.Range("F25").Value = 59
$ws.Range("G25").Value = 101
$ws.Range("H25").Value = -41.584158415841
$ws.Range("I25").Value = 363
$ws.Range("J25").Value = 499
$ws.Range("K25").Value = -27.254509018036
$ws.Range("L25").Value = -20.394736842105

# ---------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 42.857142857142
$ws.Range("F26").Value = 38
$ws.Range("G26").Value = 30
$ws.Range("H26").Value = 26.666666666666
$ws.Range("I26").Value = 126
$ws.Range("J26").Value = 128
$ws.Range("K26").Value = -1.5625
$ws.Range("L26").Value = -17.105263157894
$ws.Range("M26").Value = 50

# ---------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 75
$ws.Range("I28").Value = 29
$ws.Range("J28").Value = 22
$ws.Range("K28").Value = 31.818181818181
$ws.Range("L28").Value = 16

# ---------------------------------------------------------------------
# Row 31 - Hate Crimes
# ---------------------------------------------------------------------
Set-TextZero "C31"
$ws.Range("E31").Value = -100
$ws.Range("G31").Value = 6
$ws.Range("H31").Value = -66.666666666666
$ws.Range("J31").Value = 9
$ws.Range("K31").Value = -33.333333333333

$excel.CutCopyMode = $false
